# Auto-generated from the OOXML diff: updates crypto price/volume figures
# (and the swapped WrappedEther/Polkadot rows 12-13) in Sheet1.
#
# All target cells are strings (t="inlineStr" in the source), even when
# their text looks numeric (e.g. "1.003"), so each write forces Text
# number format first and restores the "Normal" style afterwards. This
# stops Excel from auto-coercing number-looking strings into real numbers
# while leaving the cell style index untouched (the diff shows no style
# changes).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "26.323.95"
Set-TextValue "E2" "  -5.76%  "
Set-TextValue "D3" "1.677.69"
Set-TextValue "E3" "  -3.50%  "
Set-TextValue "D4" "1.003"
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "218.29"
Set-TextValue "E5" "  -3.57%  "
Set-TextValue "D6" "0.5112"
Set-TextValue "E6" "  -11.67%  "
Set-TextValue "D7" "1.004"
Set-TextValue "E7" "  +0.14%  "
Set-TextValue "D8" "0.2669"
Set-TextValue "E8" "  -2.11%  "
Set-TextValue "D9" "0.06376"
Set-TextValue "E9" "  -3.41%  "
Set-TextValue "D10" "21.71"
Set-TextValue "E10" "  -6.66%  "
Set-TextValue "D11" "0.07375"
Set-TextValue "E11" "  -2.41%  "
Set-TextValue "B12" "WrappedEther"
Set-TextValue "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.678.68"
Set-TextValue "E12" "  -3.49%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.558"
Set-TextValue "E13" "  -3.00%  "
Set-TextValue "D14" "0.5833"
Set-TextValue "E14" "  -3.04%  "
Set-TextValue "D15" "1.896.92"
Set-TextValue "E15" "  -3.94%  "
Set-TextValue "D16" "0.000008609"
Set-TextValue "E16" "  -1.22%  "
Set-TextValue "D17" "65.35"
Set-TextValue "E17" "  -12.35%  "
Set-TextValue "D18" "26.378.74"
Set-TextValue "E18" "  -5.57%  "
Set-TextValue "D19" "4.960"
Set-TextValue "E19" "  -6.68%  "
Set-TextValue "D21" "10.87"
Set-TextValue "E21" "  -3.56%  "
Set-TextValue "D22" "190.49"
Set-TextValue "E22" "  -7.27%  "
Set-TextValue "D23" "6.227"
Set-TextValue "E23" "  -5.88%  "
Set-TextValue "D24" "1.004"
Set-TextValue "E24" "  +0.04%  "
Set-TextValue "D25" "144.16"
Set-TextValue "E25" "  -3.87%  "
Set-TextValue "D26" "7.708"
Set-TextValue "E26" "  -5.17%  "
Set-TextValue "D27" "0.1177"
Set-TextValue "E27" "  -4.32%  "
Set-TextValue "D28" "15.78"
Set-TextValue "E28" "  -2.15%  "
Set-TextValue "D29" "0.05831"
Set-TextValue "E29" "  -5.22%  "
Set-TextValue "D30" "1.284"
Set-TextValue "E30" "  -7.19%  "
Set-TextValue "D31" "1.323"
Set-TextValue "E31" "  -4.88%  "
Set-TextValue "D32" "3.539"
Set-TextValue "E32" "  -4.93%  "
Set-TextValue "D33" "3.526"
Set-TextValue "E33" "  -5.58%  "
Set-TextValue "D34" "1.657"
Set-TextValue "E34" "  -0.51%  "
Set-TextValue "D35" "1.013"
Set-TextValue "E35" "  -2.10%  "
Set-TextValue "D36" "0.6015"
Set-TextValue "E36" "  -6.21%  "
Set-TextValue "D37" "2.356"
Set-TextValue "E37" "  -2.54%  "
Set-TextValue "D38" "2.653"
Set-TextValue "E38" "  -2.35%  "
Set-TextValue "D39" "0.01621"
Set-TextValue "E39" "  -2.86%  "
Set-TextValue "D40" "6.046"
Set-TextValue "E40" "  -2.07%  "
Set-TextValue "D41" "1.084.68"
Set-TextValue "E41" "  -3.89%  "
Set-TextValue "D42" "0.8615"
Set-TextValue "E42" "  -1.42%  "
Set-TextValue "D43" "1.009"
Set-TextValue "E43" "  +0.50%  "
Set-TextValue "D44" "99.95"
Set-TextValue "E44" "  +0.31%  "
Set-TextValue "D45" "1.820.55"
Set-TextValue "E45" "  -3.61%  "
Set-TextValue "D46" "0.00000000109"
Set-TextValue "E46" "  +2.06%  "
Set-TextValue "D47" "56.04"
Set-TextValue "E47" "  -5.51%  "
Set-TextValue "D48" "1.006"
Set-TextValue "E48" "  +0.89%  "
Set-TextValue "D49" "8.117"
Set-TextValue "E49" "  -1.59%  "
Set-TextValue "D51" "0.05186"
Set-TextValue "E51" "  -3.57%  "
